$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q4" sheet right after "总计" (i.e. before the
#    current #2 sheet, "2022-Q3"). All the quarter sheets behind it simply
#    shift back by one position - their own content is untouched.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q4Sheet = $wb.Worksheets.Add($beforeSheet)
$q4Sheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Populate the new sheet: clone the layout/formatting of the (now shifted)
#    "2022-Q3" sheet - identical headers/styles - then overwrite the data row
#    with the real 2022-Q4 numbers.
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(3)
$q3Sheet.Range("A1:H2").Copy($q4Sheet.Range("A1:H2"))

$q4Sheet.Range("D2").Value = "'11.28"
$q4Sheet.Range("E2").Value = "'97.07"
$q4Sheet.Range("F2").Value = "'7.20"
$q4Sheet.Range("G2").Value = "'0.8122"
$q4Sheet.Range("H2").Value = 6

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: a new row is inserted for 2022-Q4 and
#    every later row's 日期/持有数量/持有市值 slides down one slot, while the
#    leading index column (A, 0..8) simply keeps counting.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Stretch the bordered/bold formatting down onto the brand-new row 10.
$summary.Range("A9:D9").Copy($summary.Range("A10:D10"))

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.8100000000000001

$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.73

$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.9399999999999999

$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 0.77

$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 0.65

$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 3
$summary.Range("D7").Value = 0.88

$summary.Range("B8").Value = "2021-Q2"
$summary.Range("C8").Value = 3
$summary.Range("D8").Value = 1.14

$summary.Range("B9").Value = "2021-Q1"
$summary.Range("C9").Value = 5
$summary.Range("D9").Value = 1.3

$summary.Range("A10").Value = 8
$summary.Range("B10").Value = "2020-Q4"
$summary.Range("C10").Value = 3
$summary.Range("D10").Value = 1.15

# ---------------------------------------------------------------------------
# 4. Keep the originally-selected tab ("2020-Q4", now the last sheet) active -
#    inserting the new sheet would otherwise steal the selection.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()

Write-Host "done"
